$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.227.32"
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = "'1.571.28"
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = "'211.19"
$ws.Range("E5").Value = '  +1.84%  '
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").Value = "'22.06"
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = "'0.249"
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").Value = "'1.795.58"
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("D13").Value = "'1.558.15"
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = "'3.78"
$ws.Range("E14").Value = '  +0.85%  '
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("D16").Value = "'27.175.78"
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("E18").Value = '  +1.96%  '
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = "'216.39"
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("E22").Value = '  +1.23%  '
$ws.Range("D23").Value = "'9.23"
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = "'1.94"
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").Value = "'153.91"
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = "'15.08"
$ws.Range("D28").Value = "'0.106"
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("E30").Value = '  +2.25%  '
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("D33").Value = "'1.450.42"
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("E34").Value = '  +1.62%  '
$ws.Range("D35").Value = "'1.13"
$ws.Range("E35").Value = '  +7.96%  '
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").Value = "'0.534"
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").Value = "'5.84"
$ws.Range("E40").Value = '  +2.30%  '
$ws.Range("D41").Value = "'0.811"
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").Value = "'1.706.86"
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").Value = "'86.02"
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("E49").Value = '  +3.72%  '
$ws.Range("D50").Value = "'0.0524"
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").Value = "'0.0958"
$ws.Range("E51").Value = '  +0.19%  '
